$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the text "R40" and needs to become the text "1"
# (still text, not a number) without disturbing B11's existing style/format.
#
# Assigning ($ws.Range("B11").Value = "1") directly would be auto-coerced to
# a *number* by Excel's usual "looks like a number -> store as number" rule,
# which is not what we want (the target keeps a shared-string/text cell).
# Forcing text with a leading apostrophe on B11 itself would also stamp a
# "quoted text" variant onto B11's cell style. To avoid touching B11's
# formatting at all, stage the text value "1" on a scratch cell, then copy
# only the *value* (not formats) over to B11, leaving B11's style untouched.

$scratch = $ws.Range("Z100")

$scratch.NumberFormat = "@"    # text format, so the digit isn't re-parsed as a number
$scratch.Value = "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues: value/type only, keep B11's own formatting

$scratch.Clear()
